$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 ("Resilienz der Übertragung", I18 = 75%): the "TODO" justification comment
# in column J is no longer needed -> clear it.
$ws.Range("J18").Value = ""

# Row 20 ("Anzahl Partnersysteme"): fulfillment degree revised from 70% to 66%,
# and its "TODO" justification comment in column J is cleared.
$ws.Range("I20").Value = "66%"
$ws.Range("J20").Value = ""

# Row 28 ("Alerting Error-Monitoring"): "TODO" justification comment cleared.
$ws.Range("J28").Value = ""

# Row 34 ("Alerting Metriken"): "TODO" justification comment cleared.
$ws.Range("J34").Value = ""

# Reflect the author's final scroll position / selection (view was scrolled down
# so row 21 is at the top, with M33 as the active cell).
$ws.Range("M33").Select() | Out-Null
